$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BHEL purchased: update portfolio totals (row 2) to reflect the new
# purchase/current prices, gain/loss amount, and gain/loss percentage.
$ws.Range("B2").Value = 161424.2
$ws.Range("C2").Value = 170867.7984161377
$ws.Range("D2").Value = 9443.598416137684
$ws.Range("E2").Value = 5.850175138633293
